$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("effort")

# Add the new row 40: Date, Effort [h], (no Additional Effort), Description
$ws.Range("A40").Value2 = 41225
$ws.Range("A40").NumberFormat = "ddd\ dd/mm/yyyy"

$ws.Range("B40").Value = 3.5

$ws.Range("D40").Value = "Installer creation scripts continued, new test case tc08"

# Update the selected/active cell to match the new last-used cell
$ws.Range("E40").Select()
